$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Fn1", "Itgb8", "ECs", 3, 1, 19.81532133333333, 59.445964, 0.02448901563336021, 0.02448901563336021, 1, 0.3333333333333333, 0.009519999999999999, 0.02856, 0.0009583584527718872, 0.0009583584527718872, 0.1886418590933333, 1.69777673184, 0.00002346925513229365, 0.00002346925513229365),
    @("ECs", "Fn1", "Itgb8", "FAPs", 3, 1, 19.81532133333333, 59.445964, 0.02448901563336021, 0.02448901563336021, 3, 1, 1.079146666666667, 3.23744, 0.1086354338004839, 0.1086354338004838, 21.38363796579556, 192.45274169216, 0.002660374836676918, 0.002660374836676917),
    @("ECs", "Fn1", "Itgb8", "sCs", 3, 1, 19.81532133333333, 59.445964, 0.02448901563336021, 0.02448901563336021, 3, 1, 8.844986, 26.534958, 0.8904062077467443, 0.8904062077467442, 175.2662397788347, 1577.396158009512, 0.021805171541551, 0.021805171541551),
    @("FAPs", "Fn1", "Itgb8", "ECs", 3, 1, 530.2995503333333, 1590.898651, 0.6553774102381563, 0.6553774102381563, 1, 0.3333333333333333, 0.009519999999999999, 0.02856, 0.0009583584527718872, 0.0009583584527718872, 5.048451719173332, 45.43606547256, 0.0006280864808574859, 0.0006280864808574859),
    @("FAPs", "Fn1", "Itgb8", "FAPs", 3, 1, 530.2995503333333, 1590.898651, 0.6553774102381563, 0.6553774102381563, 3, 1, 1.079146666666667, 3.23744, 0.1086354338004839, 0.1086354338004838, 572.2709920770488, 5150.43892869344, 0.07119720926425979, 0.07119720926425978),
    @("FAPs", "Fn1", "Itgb8", "sCs", 3, 1, 530.2995503333333, 1590.898651, 0.6553774102381563, 0.6553774102381563, 3, 1, 8.844986, 26.534958, 0.8904062077467443, 0.8904062077467442, 4690.492098504628, 42214.42888654166, 0.5835521144930391, 0.583552114493039),
    @("M1", "Fn1", "Itgb8", "ECs", 3, 1, 114.1551646666666, 342.4654939999999, 0.1410801048907614, 0.1410801048907614, 1, 0.3333333333333333, 0.009519999999999999, 0.02856, 0.0009583584527718872, 0.0009583584527718872, 1.086757167626666, 9.780814508639997, 0.0001352053110400057, 0.0001352053110400057),
    @("M1", "Fn1", "Itgb8", "FAPs", 3, 1, 114.1551646666666, 342.4654939999999, 0.1410801048907614, 0.1410801048907614, 3, 1, 1.079146666666667, 3.23744, 0.1086354338004839, 0.1086354338004838, 123.1901654328178, 1108.71148889536, 0.01532629839542563, 0.01532629839542563),
    @("M1", "Fn1", "Itgb8", "sCs", 3, 1, 114.1551646666666, 342.4654939999999, 0.1410801048907614, 0.1410801048907614, 3, 1, 8.844986, 26.534958, 0.8904062077467443, 0.8904062077467442, 1009.700833304361, 9087.30749973925, 0.1256186011842958, 0.1256186011842957),
    @("M2", "Fn1", "Itgb8", "ECs", 3, 1, 42.33003733333333, 126.990112, 0.05231411232645103, 0.05231411232645103, 1, 0.3333333333333333, 0.009519999999999999, 0.02856, 0.0009583584527718872, 0.0009583584527718872, 0.4029819554133333, 3.62683759872, 0.00005013567174731232, 0.00005013567174731232),
    @("M2", "Fn1", "Itgb8", "FAPs", 3, 1, 42.33003733333333, 126.990112, 0.05231411232645103, 0.05231411232645103, 3, 1, 1.079146666666667, 3.23744, 0.1086354338004839, 0.1086354338004838, 45.68031868814222, 411.1228681932801, 0.005683166286471248, 0.005683166286471247),
    @("M2", "Fn1", "Itgb8", "sCs", 3, 1, 42.33003733333333, 126.990112, 0.05231411232645103, 0.05231411232645103, 3, 1, 8.844986, 26.534958, 0.8904062077467443, 0.8904062077467442, 374.4085875928106, 3369.677288335296, 0.04658081036823247, 0.04658081036823247),
    @("Neutro", "Fn1", "Itgb8", "ECs", 3, 1, 46.22062233333333, 138.661867, 0.05712234103418551, 0.05712234103418551, 1, 0.3333333333333333, 0.009519999999999999, 0.02856, 0.0009583584527718872, 0.0009583584527718872, 0.4400203246133333, 3.96018292152, 0.00005474367837223011, 0.00005474367837223011),
    @("Neutro", "Fn1", "Itgb8", "FAPs", 3, 1, 46.22062233333333, 138.661867, 0.05712234103418551, 0.05712234103418551, 3, 1, 1.079146666666667, 3.23744, 0.1086354338004839, 0.1086354338004838, 49.87883052227556, 448.9094747004801, 0.006205510297947922, 0.006205510297947922),
    @("Neutro", "Fn1", "Itgb8", "sCs", 3, 1, 46.22062233333333, 138.661867, 0.05712234103418551, 0.05712234103418551, 3, 1, 8.844986, 26.534958, 0.8904062077467443, 0.8904062077467442, 408.8207574496207, 3679.386817046587, 0.05086208705786536, 0.05086208705786535),
    @("sCs", "Fn1", "Itgb8", "ECs", 3, 1, 56.330706, 168.992118, 0.06961701587708551, 0.06961701587708552, 1, 0.3333333333333333, 0.009519999999999999, 0.02856, 0.0009583584527718872, 0.0009583584527718872, 0.5362683211199999, 4.82641489008, 0.00006671805562255957, 0.00006671805562255958),
    @("sCs", "Fn1", "Itgb8", "FAPs", 3, 1, 56.330706, 168.992118, 0.06961701587708551, 0.06961701587708552, 3, 1, 1.079146666666667, 3.23744, 0.1086354338004839, 0.1086354338004838, 60.78909361088, 547.1018424979201, 0.007562874719702357, 0.007562874719702357),
    @("sCs", "Fn1", "Itgb8", "sCs", 3, 1, 56.330706, 168.992118, 0.06961701587708551, 0.06961701587708552, 3, 1, 8.844986, 26.534958, 0.8904062077467443, 0.8904062077467442, 498.244305940116, 4484.198753461044, 0.06198742310176059, 0.0619874231017606)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
